$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update F column "想去人数" values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8488
$ws1.Range("F5").Value = 6179
$ws1.Range("F6").Value = 536
$ws1.Range("F7").Value = 114
$ws1.Range("F10").Value = 324
$ws1.Range("F11").Value = 1141

# Sheet "全部类型" (fourth sheet) - update F column "想去人数" values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8488
$ws4.Range("F5").Value = 6179
$ws4.Range("F6").Value = 536
$ws4.Range("F7").Value = 114
$ws4.Range("F10").Value = 324
$ws4.Range("F15").Value = 1141
